# Finalizado el análisis estadístico de TP4.
#
# The workbook already has a (mostly empty) 9th worksheet named
# "460-regular 5000-nodos" with only its header row filled in (A1:D1,
# shared strings "Colores"/"Secuencial"/"Matula"/"Powell"). This script
# finishes it off the same way the previous 8 "C0xx" sheets were built:
#   1. fill in the color-frequency data rows + totals row,
#   2. add the 3-series clustered-column chart (with data table) that
#      goes with it,
#   3. add the workbook-level defined name that scopes to this sheet,
#   4. make this the active/selected sheet (mirroring the "next TP" that
#      was being analyzed becoming the front sheet).

$wb = $excel.ActiveWorkbook
$sheetName = "460-regular 5000-nodos"
$ws = $wb.Worksheets.Item($sheetName)

# --- 1. data -------------------------------------------------------------

$rows = @(
    @(96,  1,    0,     0),
    @(97,  225,  0,     0),
    @(98,  3609, 0,     10000),
    @(99,  5137, 10000, 0),
    @(100, 998,  0,     0),
    @(101, 30,   0,     0)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# totals row under the data (same pattern as every other sheet in the book)
$ws.Range("B8").Formula = "=SUM(B2:B7)"
$ws.Range("C8").Formula = "=SUM(C2:C7)"
$ws.Range("D8").Formula = "=SUM(D2:D7)"

$ws.Columns.Item(2).ColumnWidth = 11.85546875

# --- 2. chart --------------------------------------------------------------

$chartObj = $ws.ChartObjects().Add(323849, 0, 3883025, 4895850)
$chartObj.Name = "1 Gráfico"
$chart = $chartObj.Chart
$chart.ChartType = 51   # xlColumnClustered

$sSecuencial = $chart.SeriesCollection().NewSeries()
$sSecuencial.Formula = "=SERIES('" + $sheetName + "'!`$B`$1,'" + $sheetName + "'!`$A`$2:`$A`$7,'" + $sheetName + "'!`$B`$2:`$B`$7,1)"

$sMatula = $chart.SeriesCollection().NewSeries()
$sMatula.Formula = "=SERIES('" + $sheetName + "'!`$C`$1,'" + $sheetName + "'!`$A`$2:`$A`$7,'" + $sheetName + "'!`$C`$2:`$C`$7,2)"

$sPowell = $chart.SeriesCollection().NewSeries()
$sPowell.Formula = "=SERIES('" + $sheetName + "'!`$D`$1,'" + $sheetName + "'!`$A`$2:`$A`$7,'" + $sheetName + "'!`$D`$2:`$D`$7,3)"

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Distribución de Colores por Frecuencia"

$chart.GapWidth = 150

$chart.HasLegend = $false
$chart.HasDataTable = $true
$dataTable = $chart.DataTable
$dataTable.ShowLegendKey = $true
$dataTable.HasBorderHorizontal = $true
$dataTable.HasBorderVertical = $true
$dataTable.HasBorderOutline = $true

$valueAxis = $chart.Axes(2, 1)
$valueAxis.HasTitle = $true
$valueAxis.AxisTitle.Text = "Frecuencia"
$valueAxis.HasMajorGridlines = $true

$categoryAxis = $chart.Axes(1, 1)

# --- 3. defined name (sheet-scoped, like the other C0xx ranges) -----------

$ws.Names.Add("_C008_grafo_460_regular_5000_nodos", "='" + $sheetName + "'!`$A`$2:`$B`$7")

# --- 4. make this the active sheet ----------------------------------------

$ws.Activate()
$ws.Range("B1:D7").Select()
